# adding note on ggplot
# Adds a new 3rd slide ("Visualisation") with a content placeholder that
# contains a note referencing "this presentation" as a hyperlink.

$p = $ppt.ActivePresentation

# Append a new slide using the "Title and Content" layout (same layout
# slide 2 / "Introduction" already uses).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder -------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Visualisation"
$title.IndentLevel = 1
$title.ParagraphFormat.Bullet.Visible = 0

# --- Content placeholder ------------------------------------------------
$content = $s.Shapes.Item(2).TextFrame.TextRange
$fullText = "Note to selves: I think the way this presentation breaks down ggplots line by line is probably the way to go!"
$content.Text = $fullText
$content.IndentLevel = 1
$content.ParagraphFormat.Bullet.Visible = 0

# Turn the "this presentation" substring into a hyperlink.
$linkText = "this presentation"
$start = $fullText.IndexOf($linkText) + 1
$len = $linkText.Length
$sub = $content.Characters($start, $len)
$sub.ActionSettings.Item(1).Hyperlink.Address = "https://evamaerey.github.io/ggplot_flipbook/ggplot_flipbook_xaringan.html"

Write-Output ("Slides: " + $p.Slides.Count)
